$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5760
$ws.Range("J3").Value = 6132
$ws.Range("F4").Value = 1900
$ws.Range("J4").Value = 1335
$ws.Range("J5").Value = 468
$ws.Range("J6").Value = 7842
$ws.Range("F7").Value = 24091
$ws.Range("J7").Value = 21537

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 374
$ws.Range("J3").Value = 414
$ws.Range("J6").Value = 458
$ws.Range("J7").Value = 1358

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 160
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 437

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 241
$ws.Range("J3").Value = 333
$ws.Range("F4").Value = 61
$ws.Range("J6").Value = 342
$ws.Range("F7").Value = 1317
$ws.Range("J7").Value = 998

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J4").Value = 23
$ws.Range("J7").Value = 662

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 157
$ws.Range("J6").Value = 194
$ws.Range("J7").Value = 546

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J2").Value = 95
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 173
$ws.Range("J5").Value = 68
$ws.Range("J6").Value = 160
$ws.Range("J7").Value = 631
$ws.Range("J8").Value = 1358
$ws.Range("J10").Value = 149
$ws.Range("J11").Value = 334
$ws.Range("J12").Value = 43
$ws.Range("J15").Value = 238
$ws.Range("J18").Value = 179
$ws.Range("J19").Value = 632
$ws.Range("J21").Value = 62
$ws.Range("J23").Value = 202
$ws.Range("J25").Value = 107
$ws.Range("J27").Value = 130
$ws.Range("J29").Value = 1203
$ws.Range("F33").Value = 1317
$ws.Range("J33").Value = 998
$ws.Range("J36").Value = 295
$ws.Range("J37").Value = 662
$ws.Range("J42").Value = 900
$ws.Range("J43").Value = 177
$ws.Range("J44").Value = 164
$ws.Range("J45").Value = 32
$ws.Range("J46").Value = 71
$ws.Range("J47").Value = 164
$ws.Range("J48").Value = 255
$ws.Range("J52").Value = 539
$ws.Range("J54").Value = 419
$ws.Range("J56").Value = 28
$ws.Range("J57").Value = 95
$ws.Range("J63").Value = 74
$ws.Range("J65").Value = 546
$ws.Range("J67").Value = 815
$ws.Range("J71").Value = 73
$ws.Range("J73").Value = 203
$ws.Range("J74").Value = 24
$ws.Range("J75").Value = 63
$ws.Range("J76").Value = 324
$ws.Range("J79").Value = 617
$ws.Range("J83").Value = 437
$ws.Range("J85").Value = 897
$ws.Range("J86").Value = 134
$ws.Range("J88").Value = 228
$ws.Range("J89").Value = 285
$ws.Range("J94").Value = 216
$ws.Range("J99").Value = 338
$ws.Range("F101").Value = 24091
$ws.Range("J101").Value = 21537

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 203
$ws.Range("J3").Value = 310
$ws.Range("J5").Value = 24
$ws.Range("J7").Value = 815

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 82
$ws.Range("J6").Value = 202
$ws.Range("J7").Value = 419

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 363
$ws.Range("J3").Value = 420
$ws.Range("J7").Value = 1203

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J2").Value = 42
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 156
$ws.Range("J3").Value = 184
$ws.Range("J7").Value = 632

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 49
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 51
$ws.Range("J6").Value = 180
$ws.Range("J7").Value = 324

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J3").Value = 42
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 194
$ws.Range("J3").Value = 183
$ws.Range("J6").Value = 465
$ws.Range("J7").Value = 900

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 202

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 173
$ws.Range("J7").Value = 617

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 295

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 194
$ws.Range("J3").Value = 189
$ws.Range("J7").Value = 631

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J3").Value = 44
$ws.Range("J6").Value = 118
$ws.Range("J7").Value = 216

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J3").Value = 32
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 107

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 101
$ws.Range("J7").Value = 334

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 70
$ws.Range("J7").Value = 203

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J4").Value = 13
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 173

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 49
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J6").Value = 87
$ws.Range("J7").Value = 285

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J2").Value = 23
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 72
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J4").Value = 20
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 177

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J6").Value = 263
$ws.Range("J7").Value = 897

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("J3").Value = 9
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 167
$ws.Range("J7").Value = 539

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 24
